$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AllPropertyTypes")

# Add new "null" value at D5
$ws.Range("D5").Value = "null"

# Change "Assert" (A18) to "Then"
$ws.Range("A18").Value = "Then"

# Update the selected cell / active cell in the sheet view
$ws.Range("B17").Select()

# Update the workbook window position/size
$excel.Windows.Item(1).Left = 3420
$excel.Windows.Item(1).Top = 3420
$excel.Windows.Item(1).Width = 21600
$excel.Windows.Item(1).Height = 11265
